# SAM TODO.xlsx update - "updated SAM to do from last meeting - needs some sprucing up"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing status cells in column A: "Not done" -> "Future"/"Done" ---
$ws.Range("A44").Value = "Future"
$ws.Range("A47").Value = "Future"
$ws.Range("A49").Value = "Future"
$ws.Range("A53").Value = "Done"
$ws.Range("A56").Value = "Future"
$ws.Range("A61").Value = "Done"
$ws.Range("A64").Value = "Future"

# --- Append the new TODO rows (84, 86-103) picked up from the last meeting ---
# Values are written in the same order the author entered them so new shared
# strings land at the same indices as the target workbook.
$ws.Range("B86").Value = "Janine has to fix what Aron broke in the Macros"
$ws.Range("B87").Value = "Janine can sort TODO list"
$ws.Range("B88").Value = "Tower fixes  - Steve & Mike"
$ws.Range("B91").Value = "Molten salt linear fresnel"
$ws.Range("B92").Value = "Reorganize dispatch widget to UI"
$ws.Range("B93").Value = "Sample files - esp Novatec - Mike"
$ws.Range("B94").Value = "Wind - cost and scaling model - Janine"
$ws.Range("B96").Value = "Check the reports"
$ws.Range("B97").Value = "Check the macros"
$ws.Range("B98").Value = "Documentation updates from Paul"
$ws.Range("B99").Value = "Progress updates for solarpilot - Aron"
$ws.Range("C88").Value = "Steve/Mike"
$ws.Range("C89").Value = "Ty/Steve"
$ws.Range("C93").Value = "Mike"
$ws.Range("B89").Value = "Steam tower"
$ws.Range("B90").Value = "Cavity receiver"
$ws.Range("C90").Value = "Ty/Mike/Steve"
$ws.Range("B95").Value = "Review default values, financial, cost #s for PV res/com/util"
$ws.Range("B100").Value = "Subhourly simulation for physical trough"
$ws.Range("C100").Value = "Aron/Mike"
$ws.Range("B101").Value = "Check all results, and summarize for release notes"
$ws.Range("B102").Value = "Add performance adjustment factors to wind model"
$ws.Range("B84").Value = "Curtailment month by hour factors in popup widget thingy"
$ws.Range("B103").Value = "Possible registration issues"

# Remaining cells for the new rows
$ws.Range("A84").Value = "Future"
$ws.Range("C84").Value = "Aron"
$ws.Range("C86").Value = "Janine"
$ws.Range("C87").Value = "Janine"
$ws.Range("C91").Value = "Steve"
$ws.Range("C92").Value = "Steve"
$ws.Range("C94").Value = "Janine"
$ws.Range("C95").Value = "Everyone"
$ws.Range("C96").Value = "Everyone"
$ws.Range("C97").Value = "Everyone"
$ws.Range("C98").Value = "Paul"
$ws.Range("C99").Value = "Aron"
$ws.Range("C101").Value = "Everyone"
$ws.Range("C102").Value = "Janine/Aron"

# --- Apply the yellow "Future row" fill (same style as B78:B82) to the new
#     column-B cells, and plain fill-less style to the new column-A cell ---
$ws.Range("B78").Copy() | Out-Null
$ws.Range("B84").PasteSpecial(-4122) | Out-Null
$ws.Range("B86").PasteSpecial(-4122) | Out-Null
$ws.Range("B87").PasteSpecial(-4122) | Out-Null
$ws.Range("B88").PasteSpecial(-4122) | Out-Null
$ws.Range("B89").PasteSpecial(-4122) | Out-Null
$ws.Range("B90").PasteSpecial(-4122) | Out-Null
$ws.Range("B91").PasteSpecial(-4122) | Out-Null
$ws.Range("B92").PasteSpecial(-4122) | Out-Null
$ws.Range("B93").PasteSpecial(-4122) | Out-Null
$ws.Range("B94").PasteSpecial(-4122) | Out-Null
$ws.Range("B95").PasteSpecial(-4122) | Out-Null
$ws.Range("B96").PasteSpecial(-4122) | Out-Null
$ws.Range("B97").PasteSpecial(-4122) | Out-Null
$ws.Range("B98").PasteSpecial(-4122) | Out-Null
$ws.Range("B99").PasteSpecial(-4122) | Out-Null
$ws.Range("B100").PasteSpecial(-4122) | Out-Null
$ws.Range("B101").PasteSpecial(-4122) | Out-Null
$ws.Range("B102").PasteSpecial(-4122) | Out-Null
$ws.Range("B103").PasteSpecial(-4122) | Out-Null

$ws.Range("A78").Copy() | Out-Null
$ws.Range("A84").PasteSpecial(-4122) | Out-Null

# --- Update view state: scroll down and select C103 like the saved file ---
$ws.Range("C103").Select() | Out-Null
